$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.589.79"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.695.60"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.45"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3941"
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4021"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.523"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.82"
$ws.Range("E11").Value = "  +8.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08773"
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.229"
$ws.Range("E13").Value = "  +7.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.22"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001326"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.586"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.698.06"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.10"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07068"
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.870"
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.04"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.582.15"
$ws.Range("E24").Value = "  +3.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.041"
$ws.Range("E25").Value = "  +8.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.313"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.39"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.01"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.216"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.51"
$ws.Range("E30").Value = "  +3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.444"
$ws.Range("E31").Value = "  +14.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.884.67"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.105"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08533"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.246"
$ws.Range("E35").Value = "  +10.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.35"
$ws.Range("E36").Value = "  +8.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.954"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02755"
$ws.Range("E40").Value = "  +9.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09058"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7720"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7203"
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.45"
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.537"
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.211"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.355"
$ws.Range("E48").Value = "  +13.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.37"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("E51").Value = "  +3.36%  "
